$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("B5")
Write-Output $ws.Range("A1").Value()
Write-Output $ws.Range("A2").Value()
Write-Output $ws.Range("A3").Value()
